$d = $word.ActiveDocument
$anchor = $d.Paragraphs.Last

# Step 1: create all 26 new (blank) paragraphs first, chained from the last
# existing paragraph, BEFORE any numbering/list formatting is applied to any of
# them -- this avoids new paragraphs inheriting list numbering from a predecessor.
for ($i = 0; $i -lt 26; $i++) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $d.Paragraphs.Last
}
$startIndex = $d.Paragraphs.Count - 26 + 1

# Step 2: fill in text and (where needed) list numbering for each new paragraph
# paragraph 0: ilvl=None text=''
$cur = $d.Paragraphs.Item($startIndex + 0)

# paragraph 1: ilvl=0 text='3 different ligand MPNN parameters'
$cur = $d.Paragraphs.Item($startIndex + 1)
$cur.Range.Text = "3 different ligand MPNN parameters"
$cur.Range.ListFormat.ApplyNumberDefault()
$lt = $cur.Range.ListFormat.ListTemplate
$lvl1 = $lt.ListLevels.Item(1)
$lvl1.NumberStyle = 4
$lvl1.NumberFormat = "%1)"
$lvl2 = $lt.ListLevels.Item(2)
$lvl2.NumberStyle = 4
$lvl2.NumberFormat = "%2."

# paragraph 2: ilvl=1 text='Normal'
$cur = $d.Paragraphs.Item($startIndex + 2)
$cur.Range.Text = "Normal"
$cur.Range.ListFormat.ApplyNumberDefault()
$cur.Range.ListFormat.ListLevelNumber = 2

# paragraph 3: ilvl=1 text='Higher temperature: 0.25'
$cur = $d.Paragraphs.Item($startIndex + 3)
$cur.Range.Text = "Higher temperature: 0.25"
$cur.Range.ListFormat.ApplyNumberDefault()
$cur.Range.ListFormat.ListLevelNumber = 2

# paragraph 4: ilvl=1 text='Side-chain packing'
$cur = $d.Paragraphs.Item($startIndex + 4)
$cur.Range.Text = "Side-chain packing"
$cur.Range.ListFormat.ApplyNumberDefault()
$cur.Range.ListFormat.ListLevelNumber = 2

# paragraph 5: ilvl=None text=''
$cur = $d.Paragraphs.Item($startIndex + 5)

# paragraph 6: ilvl=None text='*****'
$cur = $d.Paragraphs.Item($startIndex + 6)
$cur.Range.Text = "*****"

# paragraph 7: ilvl=None text='8A ligand mpnn results, choosing by overall confidence (4/6/2024):'
$cur = $d.Paragraphs.Item($startIndex + 7)
$cur.Range.Text = "8A ligand mpnn results, choosing by overall confidence (4/6/2024):"

# paragraph 8: ilvl=None text='structure 1 –  id3: '
$cur = $d.Paragraphs.Item($startIndex + 8)
$cur.Range.Text = "structure 1 –  id3: "

# paragraph 9: ilvl=None text='a) default: id7'
$cur = $d.Paragraphs.Item($startIndex + 9)
$cur.Range.Text = "a) default: id7"

# paragraph 10: ilvl=None text='b) side chain packing: id9'
$cur = $d.Paragraphs.Item($startIndex + 10)
$cur.Range.Text = "b) side chain packing: id9"

# paragraph 11: ilvl=None text='c) higher temp: id6'
$cur = $d.Paragraphs.Item($startIndex + 11)
$cur.Range.Text = "c) higher temp: id6"

# paragraph 12: ilvl=None text=''
$cur = $d.Paragraphs.Item($startIndex + 12)

# paragraph 13: ilvl=None text='Structure 3 – id5: '
$cur = $d.Paragraphs.Item($startIndex + 13)
$cur.Range.Text = "Structure 3 – id5: "

# paragraph 14: ilvl=None text='a) default: id2'
$cur = $d.Paragraphs.Item($startIndex + 14)
$cur.Range.Text = "a) default: id2"

# paragraph 15: ilvl=None text='b) side chain packing: id7'
$cur = $d.Paragraphs.Item($startIndex + 15)
$cur.Range.Text = "b) side chain packing: id7"

# paragraph 16: ilvl=None text='c) higer temp: id4'
$cur = $d.Paragraphs.Item($startIndex + 16)
$cur.Range.Text = "c) higer temp: id4"

# paragraph 17: ilvl=None text=''
$cur = $d.Paragraphs.Item($startIndex + 17)

# paragraph 18: ilvl=None text='MPNN only - id7:'
$cur = $d.Paragraphs.Item($startIndex + 18)
$cur.Range.Text = "MPNN only - id7:"

# paragraph 19: ilvl=None text='a) default: id2'
$cur = $d.Paragraphs.Item($startIndex + 19)
$cur.Range.Text = "a) default: id2"

# paragraph 20: ilvl=None text='b) side chain packing: id7'
$cur = $d.Paragraphs.Item($startIndex + 20)
$cur.Range.Text = "b) side chain packing: id7"

# paragraph 21: ilvl=None text='c) higer temp: id4'
$cur = $d.Paragraphs.Item($startIndex + 21)
$cur.Range.Text = "c) higer temp: id4"

# paragraph 22: ilvl=None text=''
$cur = $d.Paragraphs.Item($startIndex + 22)

# paragraph 23: ilvl=None text=''
$cur = $d.Paragraphs.Item($startIndex + 23)

# paragraph 24: ilvl=None text='*****'
$cur = $d.Paragraphs.Item($startIndex + 24)
$cur.Range.Text = "*****"

# paragraph 25: ilvl=None text=''
$cur = $d.Paragraphs.Item($startIndex + 25)

Write-Output "paragraphs=$($d.Paragraphs.Count)"
